$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O21").Value = 0.0
$ws.Range("S21").Value = 130.0

$ws.Range("R22").Value = 1500.0
$ws.Range("S22").Value = 1500.0

$ws.Range("T23").Value = 1500.0

$ws.Range("L24").Value = 4958.0
$ws.Range("M24").Value = 498.0

$ws.Range("O25").Value = 5456.0
$ws.Range("Q25").Value = 546.0
